# Applies the "Final plots and data for the gamma ray lab" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Al" worksheet (2nd sheet) - view only changed: scroll/selection moved.
# ---------------------------------------------------------------------------
$wsAl = $wb.Worksheets.Item(2)
$wsAl.Activate()
$wsAl.Range("A7").Select()
$wsAl.Range("C11").Select()

# ---------------------------------------------------------------------------
# 2. "energy-v-murho" worksheet (3rd sheet) - updated measured B values and
#    formulas/results in column C, plus several newly filled-in rows.
# ---------------------------------------------------------------------------
$wsEvm = $wb.Worksheets.Item(3)
$wsEvm.Activate()

$wsEvm.Range("B2").Value = 0.4029
$wsEvm.Range("C2").Formula = "=(0.068)*B2"

$wsEvm.Range("B3").Value = 0.3378
$wsEvm.Range("C3").Formula = "=B3*(0.032)"

$wsEvm.Range("B4").Value = 0.2492
$wsEvm.Range("C4").Formula = "=B4*(0.018)"

$wsEvm.Range("B5").Value = 0.2032
$wsEvm.Range("C5").Formula = "=B5*(0.043)"

$wsEvm.Range("B6").Value = 0.0982
$wsEvm.Range("C6").Formula = "=B6*(0.0109)"

$wsEvm.Range("B7").Value = 0.0704
$wsEvm.Range("C7").Formula = "=B7*(0.054)"

$wsEvm.Range("B8").Value = 0.0507
$wsEvm.Range("C8").Formula = "=B8*(0.044)"

$wsEvm.Range("B12").Value = 0.2226
$wsEvm.Range("C12").Formula = "=0.252*B12"

$wsEvm.Range("B13").Value = 0.0952
$wsEvm.Range("C13").Formula = "=0.1073*B13"

$wsEvm.Range("B14").Value = 0.1547
$wsEvm.Range("C14").Formula = "=0.1115*B14"

$wsEvm.Range("B15").Value = 0.1435
$wsEvm.Range("C15").Formula = "=0.1091*B15"

$wsEvm.Range("B16").Value = 0.0963
$wsEvm.Range("C16").Formula = "=0.1091*B16"

$wsEvm.Range("B17").Value = 0.0757
$wsEvm.Range("C17").Formula = "=B17*0.0844"

$wsEvm.Range("B18").Value = 0.106
$wsEvm.Range("C18").Formula = "=B18*0.2617"

$wsEvm.Range("B19").Value = 0.1105
$wsEvm.Range("C19").Formula = "=B19*0.2155"

$wsEvm.Range("I17").Select()

# ---------------------------------------------------------------------------
# 3. Brand new "Sheet1" worksheet added at the end of the workbook, holding
#    isotope reference data, and made the active tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Sheet1"

$wsNew.Range("A2").Value = "Istope: Pb, Z: 82"
$wsNew.Range("B3").Value = "g/mol"
$wsNew.Range("B4").Value = 207.2

$wsNew.Range("A6").Value = "Istope: Al, Z: 13"
$wsNew.Range("B7").Value = "g/mol"
$wsNew.Range("B8").Value = 26.3

$wsNew.Activate()
$wsNew.Range("A7").Select()
$wsNew.Range("B9").Select()
